# "Generate Report for Handback"
#
# For the zh-cn and de-de locale sheets, row 5 (the 2c43aa88-... handback
# entry) gets a freshly generated handback report:
#   - I5 (Latest Target File): a new hyperlink to the .md handback file that
#     was actually received (not yet the "latest" revision)
#   - J5 (Latest Handback File): the generated .xlf handback file name
#   - K5 (Latest Handback DateTime): the handback timestamp
#   - P5 (Error Detail): a warning that the handback is behind the latest
#     revision of the source .md file
#
# Also widens columns I, J and P (Latest Target File / Latest Handback File /
# Error Detail) to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$currentHandbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f977141640c4d45dafb7e6cc4aea7bc120efa21a/e2e/2c43aa88-3c4b-4cf3-9cf1-4a04075df26c.md"
$latestHandbackUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/970eee774b6128a9f208bcfc9cf2ed004c32a2e5/e2e/2c43aa88-3c4b-4cf3-9cf1-4a04075df26c.md"
$errorDetail = "The version of handback file is not the latest, current: $currentHandbackUrl, latest: $latestHandbackUrl."
$mdDisplay = "2c43aa88-3c4b-4cf3-9cf1-4a04075df26c.md"

function Update-LocaleSheet($SheetName, $HandbackXlf, $HandbackDateTime) {

    $ws = $wb.Worksheets.Item($SheetName)

    # Widen the columns that now hold longer generated report text.
    $ws.Columns.Item(9).ColumnWidth = 39.17   # I: Latest Target File
    $ws.Columns.Item(10).ColumnWidth = 39.17  # J: Latest Handback File
    $ws.Columns.Item(16).ColumnWidth = 39.17  # P: Error Detail

    # I5: Latest Target File - link to the handback markdown file as received.
    $ws.Range("I5").Value = $mdDisplay
    $ws.Hyperlinks.Add($ws.Range("I5"), $currentHandbackUrl, "", "", $mdDisplay) | Out-Null

    # J5: Latest Handback File - the generated xlf handback file name.
    $ws.Range("J5").Value = $HandbackXlf

    # K5: Latest Handback DateTime.
    $ws.Range("K5").Value = $HandbackDateTime

    # P5: Error Detail - handback is behind the latest source revision.
    $ws.Range("P5").Value = $errorDetail
}

Update-LocaleSheet "zh-cn" "2c43aa88-3c4b-4cf3-9cf1-4a04075df26c.befb64ffe48bef804b6d13bf660c26a6ef409012.zh-cn.xlf" "2016-10-27 07:25:58"
Update-LocaleSheet "de-de" "2c43aa88-3c4b-4cf3-9cf1-4a04075df26c.befb64ffe48bef804b6d13bf660c26a6ef409012.de-de.xlf" "2016-10-27 07:26:15"
